$d = $word.ActiveDocument

# --- Paragraph 1 (title line with manual line break): targeted text substitutions ---
$tOld0 = "🚀המאמר היומי של מייק -08.10.24: ⚡️🚀"
$tNew0 = "🚀המאמר היומי של מייק -06.10.24: ⚡️🚀"
$tFound0 = $d.Content.Find.Execute($tOld0, $true, $false, $false, $false, $false, $true, 1, $false, $tNew0, 2)
Write-Output "Title replace 0: $tFound0"

$tOld1 = "CONTEXTUAL DOCUMENT EMBEDDINGS"
$tNew1 = "CONTRASTIVE LOCALIZED LANGUAGE-IMAGE PRE-TRAINING"
$tFound1 = $d.Content.Find.Execute($tOld1, $true, $false, $false, $false, $false, $true, 1, $false, $tNew1, 2)
Write-Output "Title replace 1: $tFound1"

# --- Body paragraphs 2..10: replace full paragraph text via Delete+InsertAfter so
#     xml:space + literal apostrophes are computed fresh (avoids Find/Replace smart-quoting) ---
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
[void]$r2.MoveEnd(1, -1)
$r2.Delete()
$p2b = $d.Paragraphs.Item(2)
$p2b.Range.InsertAfter("ממשיכים הפסקה בסקירות על מודלי שפה ועוברים לסקירות על מודלים מולטימודליים (שפה ותמונות). טוב, הפסקה למחצה. אתם בטח זוכרים את המודל שנקרא CLIP שעשה הרבה רעש לפני כמה שנים. ")
Write-Output "Paragraph 2 replaced"

$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
[void]$r3.MoveEnd(1, -1)
$r3.Delete()
$p3b = $d.Paragraphs.Item(3)
$p3b.Range.InsertAfter("CLIP הוא אחד המודלים מולטימודליים הראשוניים שהצליח לייצר אמבדינגס חזקים ומיושרים (aligned) של טקסט ושל תמונות. מיושרים הכוונה של הייצוגים של תמונה וטקסט שמתאר את תוכנה קרובים אחד לשני בזמן שהייצוגים של תמונה וטקסט לא מתאימים רחוקים אחד מהשני (במקרה הזה ביחס למרחק קוסיין ביניהם).")
Write-Output "Paragraph 3 replaced"

$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
[void]$r4.MoveEnd(1, -1)
$r4.Delete()
$p4b = $d.Paragraphs.Item(4)
$p4b.Range.InsertAfter("המודל הזה אומן על דאטהסט ענק של תמונות והכותרות שלהם (או טאגים) מהאינטרנט כאשר אימנו אותו תוך שימוש בטכניקה למידה ניגודית (contrastive learning או CL). בגדול מאוד טכניקות CL מאומנות להפיק ייצוג סמנטי מדאטה (מסוגים שונים) כאשר המטרה היא לקרב את הייצוגים (אמבדינגס) של פיסות דאטה קרובות (או חיוביות) ולהרחיק ייצוגים של פיסות דאטה לא דומות (שליליות). במקרה של CLIP פיסות דאטה חיוביות הם הייצוגים של תמונה והכותרת שלה ואילו הזוגות השליליים בנויים מכותבות ותמונות שנבחרו באקראי.")
Write-Output "Paragraph 4 replaced"

$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
[void]$r5.MoveEnd(1, -1)
$r5.Delete()
$p5b = $d.Paragraphs.Item(5)
$p5b.Range.InsertAfter("המאמר שנסקור אחד כאמור משכלל את CLIP על ידי הקניה של יכולות לוקליזציה לייצוג. הכוונה כאן שהמחברים מאמנים ייצוגים של תמונה ושל טקסט באופן כזה שבהינתן ייצוג התמונה I וייצוג התיאור של פאץ' ב I המכיל אובייקט מסוים יהיה ניתן להפיק ב״קלות״ את מיקום האובייקט בתמונה. ")
Write-Output "Paragraph 5 replaced"

$p6 = $d.Paragraphs.Item(6)
$r6 = $p6.Range
[void]$r6.MoveEnd(1, -1)
$r6.Delete()
$p6b = $d.Paragraphs.Item(6)
$p6b.Range.InsertAfter("במילים פשוטות נניח שיש לנו אריה עומד ושואג בתמונה הנמצא ב-bounding box (המוגדר על ידי רביעיה של קואורדינטות שלו בתמונה) המסומן ב- B. המחברים מאמנים רשת אנקודר לתמונות f_I רשת אנקודר לטקסט f_T כך שייצוג התמונה R_I ייצוג ״אריה עומד ושואג״ R_T, המופקים על ידי שני האנקודר האלו (בהתאמה) כך שרשת רדודה יחסית (נקראת prompter במאמר), המקבלת אותם, תוכל לחזות את מיקום האריה B בתמונה. דרך אגב המיקום כאן לא חייב להיות מתואר על ידי bounding box אלא יכול להיות מוגדר (בערך) על ידי כמה ניקודת, תיאור כללי (נגיד חיה, בלי להזכיר שזה אריה) ובעוד צורות.")
Write-Output "Paragraph 6 replaced"

$p7 = $d.Paragraphs.Item(7)
$r7 = $p7.Range
[void]$r7.MoveEnd(1, -1)
$r7.Delete()
$p7b = $d.Paragraphs.Item(7)
$p7b.Range.InsertAfter("האימון נעשה כמו בלמידה הניגודית כמו ב-CLIP המקורי. אבל בנוסף ללוס הרגיל שלו יש כאן עוד לוס ניגודי המקרב את ייצוגים של כותרת הפאץ' בתמונה לייצוג המופק על Prompter מייצוג התמונה ומהמתאר של הפאץ' (נגיד BB) ומרחיק את הייצוגים האלו לפאצ'ים שונים. כמובן שה-Prompter גם מאומן תוך כדי,")
Write-Output "Paragraph 7 replaced"

$p8 = $d.Paragraphs.Item(8)
$r8 = $p8.Range
[void]$r8.MoveEnd(1, -1)
$r8.Delete()
$p8b = $d.Paragraphs.Item(8)
$p8b.Range.InsertAfter("המאמר משתמש במודלים מאומנים למטרת זיהוי אובייקטים בתמונה (OWLv2) ובמודלים מאומנים אחרים (VeCap) למתן כותרות לפאצ'ים האלו. ")
Write-Output "Paragraph 8 replaced"

$p9 = $d.Paragraphs.Item(9)
$r9 = $p9.Range
[void]$r9.MoveEnd(1, -1)
$r9.Delete()
$p9b = $d.Paragraphs.Item(9)
$p9b.Range.InsertAfter("מאמר די חמוד וקליל…")
Write-Output "Paragraph 9 replaced"

$p10 = $d.Paragraphs.Item(10)
$r10 = $p10.Range
[void]$r10.MoveEnd(1, -1)
$r10.Delete()
$p10b = $d.Paragraphs.Item(10)
$p10b.Range.InsertAfter("https://arxiv.org/pdf/2410.02746")
Write-Output "Paragraph 10 replaced"

# --- Remove the now-orphaned final paragraph that held the old arxiv link ---
$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
Write-Output "Last paragraph text before delete: $($lastP.Range.Text)"
$lastP.Range.Delete()
Write-Output "Paragraph count after delete: $($d.Paragraphs.Count)"